$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet handles
# ---------------------------------------------------------------------------
$wsMaster = $wb.Worksheets.Item(1)          # "Master" - untouched
$wsModelo = $wb.Worksheets.Item(2)          # "SchoolarModeloControlador"
$wsCampos = $wb.Worksheets.Item(3)          # "ScholarCampos"

$nl = [char]10

# ===========================================================================
# 1. SchoolarModeloControlador (sheet2): add 5 new tables (permiso, usuarioRol,
#    modulo, rol, permisorol) to the model/controller generator list.
# ===========================================================================

$wsModelo.Range("C10").Value = "permiso"
$wsModelo.Range("C11").Value = "usuarioRol"
$wsModelo.Range("C12").Value = "modulo"
$wsModelo.Range("C13").Value = "rol"
$wsModelo.Range("C14").Value = "permisorol"

# Extend the two shared formulas (model / controller generators) down through
# row 14, matching the existing pattern used for rows 4-9.
$wsModelo.Range("G4:G14").Formula = '="php artisan make:model "&PROPER(C4)&" --migration"'
$wsModelo.Range("H4:H14").Formula = '="php artisan make:controller "&PROPER(C4)&"Controller --resource"'

# Column G needs to widen to fit the newly added longer text.
$wsModelo.Columns("G").AutoFit()

# The active cell / tab moves on to the ScholarCampos sheet, so this sheet is
# no longer the tab-selected one; selection parks on C11.
$wsModelo.Range("C11").Select()

# ===========================================================================
# 2. ScholarCampos (sheet3): brand-new migration-field table describing the
#    "permiso"/"rol" (permissions/roles) schema, plus the $table-> generator
#    formulas in columns L/M.
# ===========================================================================

# Shift everything down one row (row 2 header -> row 3) to make room for a
# blank first row that holds a single highlighted helper cell (L1).
$wsCampos.Rows(1).Insert()

# --- header row (was row 2, now row 3) : new columns F/I/J/L -----------------
$wsCampos.Range("F3").Value = "porte"
$wsCampos.Range("I3").Value = "unsigned"
$wsCampos.Range("J3").Value = "foreign plural"
$wsCampos.Range("L3").Value = "gen"

# --- data rows --------------------------------------------------------------
# modulo
$wsCampos.Range("C4").Value = "modulo"
$wsCampos.Range("D4").Value = "modulo"
$wsCampos.Range("E4").Value = "string"
$wsCampos.Range("F4").Value = 100
$wsCampos.Range("H4").Value = 0

$wsCampos.Range("C5").Value = "modulo"
$wsCampos.Range("D5").Value = "activo"
$wsCampos.Range("E5").Value = "boolean"

# permiso
$wsCampos.Range("C6").Value = "permiso"
$wsCampos.Range("D6").Value = "permiso"
$wsCampos.Range("E6").Value = "string"
$wsCampos.Range("F6").Value = 100

$wsCampos.Range("C7").Value = "permiso"
$wsCampos.Range("D7").Value = "modulo_id"
$wsCampos.Range("E7").Value = "integer"
$wsCampos.Range("I7").Value = 1
$wsCampos.Range("J7").Value = "modulos"

$wsCampos.Range("C8").Value = "permiso"
$wsCampos.Range("D8").Value = "descripcion"
$wsCampos.Range("E8").Value = "string"
$wsCampos.Range("F8").Value = 150

# permisorol
$wsCampos.Range("C11").Value = "permisorol"
$wsCampos.Range("D11").Value = "permiso_id"
$wsCampos.Range("E11").Value = "integer"
$wsCampos.Range("I11").Value = 1
$wsCampos.Range("J11").Value = "permisos"

$wsCampos.Range("C12").Value = "permisorol"
$wsCampos.Range("D12").Value = "rol_id"
$wsCampos.Range("E12").Value = "integer"
$wsCampos.Range("I12").Value = 1
$wsCampos.Range("J12").Value = "rols"

# rol
$wsCampos.Range("C13").Value = "rol"
$wsCampos.Range("D13").Value = "rol"
$wsCampos.Range("E13").Value = "string"
$wsCampos.Range("F13").Value = 100

$wsCampos.Range("C14").Value = "rol"
$wsCampos.Range("D14").Value = "descripcion"
$wsCampos.Range("E14").Value = "string"
$wsCampos.Range("F14").Value = 150

# usuarioRol
$wsCampos.Range("C15").Value = "usuarioRol"
$wsCampos.Range("D15").Value = "rol_id"
$wsCampos.Range("E15").Value = "integer"
$wsCampos.Range("I15").Value = 1
$wsCampos.Range("J15").Value = "rols"

$wsCampos.Range("C16").Value = "usuarioRol"
$wsCampos.Range("D16").Value = "user_id"
$wsCampos.Range("E16").Value = "integer"
$wsCampos.Range("I16").Value = 1
$wsCampos.Range("J16").Value = "users"

# --- generator formulas in L (migration line) and M (fillable array entry) -
# Column L builds the "$table->column(...)" migration line; written first
# across L4:L9, then the fill handle is dragged further to cover L9:L23.
$formulaL_base = '="$table->"&E4&"(''"&D4&"''"&IF(LEN(F4)>0,","&F4&"","")&")"&IF(H4>0,"->nullable()","")&IF(LEN(G4)>0,"->default(''"&G4&"'')","")&IF(LEN(I4)>0,"->unsigned()","")&";"&IF(LEN(J4)>0,"' + $nl + '   $table->foreign(''"&D4&"'')->references(''id'')->on(''"&J4&"'');","")'
$wsCampos.Range("L4:L9").Formula = $formulaL_base

# Column M builds the "'field'," fillable-array snippet; filled once from
# M4 all the way down to M23.
$formulaM_base = '="''"&D4&"'',"'
$wsCampos.Range("M4:M23").Formula = $formulaM_base

# The fill handle on L gets dragged again from L9 down to L23, extending the
# generator formula (now anchored on row 10) across the rest of the table.
$formulaL_ext = '="$table->"&E10&"(''"&D10&"''"&IF(LEN(F10)>0,","&F10&"","")&")"&IF(H10>0,"->nullable()","")&IF(LEN(G10)>0,"->default(''"&G10&"'')","")&IF(LEN(I10)>0,"->unsigned()","")&";"&IF(LEN(J10)>0,"' + $nl + '   $table->foreign(''"&D10&"'')->references(''id'')->on(''"&J10&"'');","")'
$wsCampos.Range("L9:L23").Formula = $formulaL_ext

# --- formatting --------------------------------------------------------------
# L1 is a highlighted scratch cell (yellow fill + wrap text).
$wsCampos.Range("L1").Interior.Color = 65535
$wsCampos.Range("L1").WrapText = $true

# Column widths for the "bestFit" columns (C, D, L) widen to fit new content.
$wsCampos.Columns("C").AutoFit()
$wsCampos.Columns("D").AutoFit()
$wsCampos.Columns("L").AutoFit()

# Portrait page orientation.
$wsCampos.PageSetup.Orientation = 1

# This sheet becomes the active / tab-selected one, with the cursor left on I6.
$wsCampos.Activate()
$wsCampos.Range("I6").Select()
